$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The columns H (mes-nombre), J (ocupacion-1-digito-descripcion), L (ano-cno) and
# M (sexo) are re-curated from "dimension" concepts that needed an external mapping
# file to plain "measure" concepts that don't.
#
# Row 2 holds the iaest concept URI for each column.
$ws.Range("H2").Value = "iaest-measure:mes-nombre"
$ws.Range("J2").Value = "iaest-measure:ocupacion-1-digito-descripcion"
$ws.Range("L2").Value = "iaest-measure:ano-cno"
$ws.Range("M2").Value = "iaest-measure:sexo"

# Row 3 holds whether the column is a "dim" or a "medida".
$ws.Range("H3").Value = "medida"
$ws.Range("J3").Value = "medida"
$ws.Range("L3").Value = "medida"
$ws.Range("M3").Value = "medida"

# Row 4 holds the XML datatype of the column ("skos:Concept" for curated
# dimensions, "xsd:int" for plain measures).
$ws.Range("H4").Value = "xsd:int"
$ws.Range("J4").Value = "xsd:int"
$ws.Range("L4").Value = "xsd:int"
$ws.Range("M4").Value = "xsd:int"

# Row 5 used to reference the per-dimension mapping workbook; measures don't
# need one any more, so those cells are removed outright.
$ws.Range("H5").Clear()
$ws.Range("J5").Clear()
$ws.Range("L5").Clear()
$ws.Range("M5").Clear()
